$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.84
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 4.1
$ws.Range("L2").Value = 1.39
$ws.Range("O2").Value = 1.32
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 3.35
$ws.Range("T2").Value = 1.8
$ws.Range("X2").Value = 17.5
$ws.Range("Y2").Value = 19.5
$ws.Range("AA2").Value = 130
$ws.Range("AB2").Value = 11
$ws.Range("AC2").Value = 980
$ws.Range("AD2").Value = 1000
$ws.Range("AG2").Value = 12
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 70
$ws.Range("AM2").Value = 130
$ws.Range("AO2").Value = 70
# Row 3
$ws.Range("G3").Value = 2.8
$ws.Range("H3").Value = 3.05
$ws.Range("I3").Value = 3.8
$ws.Range("J3").Value = 2.84
$ws.Range("K3").Value = 3.55
$ws.Range("L3").Value = 1.43
$ws.Range("M3").Value = 1.09
$ws.Range("N3").Value = 2.9
$ws.Range("O3").Value = 1.41
$ws.Range("Q3").Value = 2.22
$ws.Range("S3").Value = 3.85
$ws.Range("V3").Value = 1.35
# Row 4
$ws.Range("F4").Value = 1.16
$ws.Range("G4").Value = 1.22
$ws.Range("H4").Value = 21
$ws.Range("I4").Value = 55
$ws.Range("J4").Value = 7.2
$ws.Range("K4").Value = 11
$ws.Range("L4").Value = 1.26
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 4.4
$ws.Range("R4").Value = 1.46
$ws.Range("S4").Value = 2.7
$ws.Range("T4").Value = 2.8
$ws.Range("U4").Value = 1.43
$ws.Range("W4").Value = 5.4
$ws.Range("X4").Value = 980
$ws.Range("Y4").Value = 80
$ws.Range("AC4").Value = 980
$ws.Range("AD4").Value = 130
$ws.Range("AF4").Value = 7.6
$ws.Range("AH4").Value = 80
$ws.Range("AK4").Value = 980
$ws.Range("AL4").Value = 90
$ws.Range("AN4").Value = 5.1
# Row 5
$ws.Range("G5").Value = 1.86
$ws.Range("H5").Value = 4.5
$ws.Range("I5").Value = 5.3
$ws.Range("J5").Value = 3.85
$ws.Range("K5").Value = 4.5
$ws.Range("L5").Value = 1.25
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 4.7
$ws.Range("O5").Value = 1.21
$ws.Range("P5").Value = 2.24
$ws.Range("Q5").Value = 1.64
$ws.Range("R5").Value = 1.52
$ws.Range("S5").Value = 2.52
$ws.Range("T5").Value = 1.65
$ws.Range("U5").Value = 2.26
$ws.Range("V5").Value = 1.24
$ws.Range("W5").Value = 2.16
$ws.Range("X5").Value = 26
$ws.Range("Y5").Value = 26
$ws.Range("Z5").Value = 48
$ws.Range("AA5").Value = 130
$ws.Range("AB5").Value = 14
$ws.Range("AC5").Value = 11.5
$ws.Range("AD5").Value = 24
$ws.Range("AE5").Value = 70
$ws.Range("AF5").Value = 15
$ws.Range("AG5").Value = 13
$ws.Range("AH5").Value = 21
$ws.Range("AI5").Value = 65
$ws.Range("AJ5").Value = 24
$ws.Range("AK5").Value = 21
$ws.Range("AL5").Value = 32
$ws.Range("AM5").Value = 90
$ws.Range("AN5").Value = 10.5
$ws.Range("AO5").Value = 60
# Row 6
$ws.Range("S6").Value = 5.6
$ws.Range("Y6").Value = 12
$ws.Range("AH6").Value = 24
# Row 7
$ws.Range("F7").Value = 2.02
$ws.Range("G7").Value = 2.26
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 3.75
$ws.Range("J7").Value = 3.75
$ws.Range("K7").Value = 4.5
$ws.Range("L7").Value = 1.25
$ws.Range("N7").Value = 5
$ws.Range("O7").Value = 1.19
$ws.Range("P7").Value = 2.38
$ws.Range("Q7").Value = 1.58
$ws.Range("R7").Value = 1.55
$ws.Range("S7").Value = 2.44
$ws.Range("T7").Value = 1.54
$ws.Range("U7").Value = 2.44
$ws.Range("V7").Value = 1.36
$ws.Range("W7").Value = 1.79
$ws.Range("X7").Value = 29
$ws.Range("Y7").Value = 23
$ws.Range("Z7").Value = 34
$ws.Range("AA7").Value = 75
$ws.Range("AB7").Value = 16.5
$ws.Range("AC7").Value = 10.5
$ws.Range("AD7").Value = 18.5
$ws.Range("AE7").Value = 42
$ws.Range("AF7").Value = 19.5
$ws.Range("AG7").Value = 13.5
$ws.Range("AH7").Value = 18.5
$ws.Range("AI7").Value = 980
$ws.Range("AJ7").Value = 32
$ws.Range("AK7").Value = 24
$ws.Range("AL7").Value = 34
$ws.Range("AM7").Value = 70
$ws.Range("AN7").Value = 13
$ws.Range("AO7").Value = 30
# Row 8
$ws.Range("F8").Value = 2.32
$ws.Range("G8").Value = 2.72
$ws.Range("H8").Value = 2.64
$ws.Range("I8").Value = 3.2
$ws.Range("N8").Value = 2.66
$ws.Range("P8").Value = 2.64
$ws.Range("Q8").Value = 1.33
$ws.Range("R8").Value = 1.83
$ws.Range("S8").Value = 1.83
$ws.Range("T8").Value = 1.29
$ws.Range("V8").Value = 1.46
$ws.Range("W8").Value = 1.58
$ws.Range("AE8").Value = 980
# Row 9
$ws.Range("G9").Value = 2.32
$ws.Range("I9").Value = 3.35
$ws.Range("J9").Value = 4.1
$ws.Range("K9").Value = 4.9
$ws.Range("L9").Value = 1.2
$ws.Range("N9").Value = 6
$ws.Range("O9").Value = 1.13
$ws.Range("Q9").Value = 1.35
$ws.Range("R9").Value = 1.82
$ws.Range("S9").Value = 1.98
$ws.Range("T9").Value = 1.42
$ws.Range("U9").Value = 2.92
$ws.Range("V9").Value = 1.42
$ws.Range("W9").Value = 1.75
$ws.Range("X9").Value = 980
$ws.Range("Y9").Value = 980
$ws.Range("Z9").Value = 980
$ws.Range("AA9").Value = 60
$ws.Range("AB9").Value = 980
$ws.Range("AD9").Value = 980
$ws.Range("AE9").Value = 980
$ws.Range("AF9").Value = 980
$ws.Range("AH9").Value = 1000
$ws.Range("AI9").Value = 980
$ws.Range("AJ9").Value = 980
$ws.Range("AK9").Value = 980
$ws.Range("AL9").Value = 980
$ws.Range("AM9").Value = 55
$ws.Range("AN9").Value = 10
$ws.Range("AO9").Value = 17.5
# Row 10
$ws.Range("G10").Value = 1.96
$ws.Range("I10").Value = 4.4
$ws.Range("K10").Value = 5.1
$ws.Range("N10").Value = 1.02
$ws.Range("O10").Value = 1.13
$ws.Range("P10").Value = 2.54
$ws.Range("Q10").Value = 1.13
$ws.Range("R10").Value = 1.68
$ws.Range("S10").Value = 1.91
$ws.Range("T10").Value = 1.01
$ws.Range("U10").Value = 1.01
$ws.Range("V10").Value = 1.29
$ws.Range("X10").Value = 1000
$ws.Range("Y10").Value = 1000
$ws.Range("Z10").Value = 1000
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 1000
$ws.Range("AC10").Value = 1000
$ws.Range("AD10").Value = 1000
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 1000
$ws.Range("AG10").Value = 1000
$ws.Range("AH10").Value = 1000
$ws.Range("AI10").Value = 1000
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 1000
$ws.Range("AL10").Value = 1000
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000
$ws.Range("AO10").Value = 1000
# Row 11
$ws.Range("G11").Value = 2.24
$ws.Range("I11").Value = 4.5
$ws.Range("J11").Value = 3.15
$ws.Range("L11").Value = 1.37
$ws.Range("N11").Value = 3.2
$ws.Range("P11").Value = 1.74
$ws.Range("Q11").Value = 2.06
$ws.Range("R11").Value = 1.28
$ws.Range("S11").Value = 3.75
$ws.Range("U11").Value = 1.97
$ws.Range("V11").Value = 1.29
$ws.Range("W11").Value = 1.81
$ws.Range("X11").Value = 15
$ws.Range("Y11").Value = 16.5
$ws.Range("AA11").Value = 110
$ws.Range("AB11").Value = 10.5
$ws.Range("AC11").Value = 8.199999999999999
$ws.Range("AD11").Value = 21
$ws.Range("AE11").Value = 70
$ws.Range("AF11").Value = 15.5
$ws.Range("AG11").Value = 13.5
$ws.Range("AI11").Value = 80
$ws.Range("AM11").Value = 140
$ws.Range("AO11").Value = 80
# Row 12
$ws.Range("F12").Value = 3.05
$ws.Range("G12").Value = 3.25
$ws.Range("H12").Value = 2.76
$ws.Range("J12").Value = 3
$ws.Range("K12").Value = 3.05
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 2.82
$ws.Range("O12").Value = 1.47
$ws.Range("P12").Value = 1.61
$ws.Range("Q12").Value = 2.36
$ws.Range("R12").Value = 1.22
$ws.Range("S12").Value = 4.7
$ws.Range("T12").Value = 1.95
$ws.Range("U12").Value = 1.89
$ws.Range("V12").Value = 1.53
$ws.Range("W12").Value = 1.45
$ws.Range("X12").Value = 12
$ws.Range("Y12").Value = 10.5
$ws.Range("Z12").Value = 20
$ws.Range("AA12").Value = 60
$ws.Range("AB12").Value = 10.5
$ws.Range("AC12").Value = 7.6
$ws.Range("AD12").Value = 15.5
$ws.Range("AE12").Value = 40
$ws.Range("AF12").Value = 21
$ws.Range("AG12").Value = 16.5
$ws.Range("AH12").Value = 22
$ws.Range("AI12").Value = 70
$ws.Range("AJ12").Value = 60
$ws.Range("AK12").Value = 55
$ws.Range("AL12").Value = 75
$ws.Range("AM12").Value = 180
$ws.Range("AN12").Value = 60
$ws.Range("AO12").Value = 50
# Row 13
$ws.Range("G13").Value = 1.28
$ws.Range("I13").Value = 15.5
$ws.Range("J13").Value = 6.8
$ws.Range("K13").Value = 7
$ws.Range("L13").Value = 1.32
$ws.Range("Q13").Value = 1.7
$ws.Range("V13").Value = 1.07
$ws.Range("W13").Value = 4.6
$ws.Range("AA13").Value = 900
$ws.Range("AM13").Value = 260
$ws.Range("AO13").Value = 450
# Row 14
$ws.Range("F14").Value = 2.84
$ws.Range("H14").Value = 2.82
$ws.Range("I14").Value = 2.86
$ws.Range("K14").Value = 3.35
$ws.Range("M14").Value = 1.09
$ws.Range("O14").Value = 1.39
$ws.Range("P14").Value = 1.8
$ws.Range("T14").Value = 1.87
$ws.Range("AB14").Value = 10.5
$ws.Range("AH14").Value = 18
